$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("C1").Value = "accepted"
$ws.Range("D1").Value = "rejected"

# Update existing row 2
$ws.Range("A2").Value = 76.55989583333334
$ws.Range("B2").Value = 2.0
$ws.Range("C2").Value = 435.0
# D2 stays empty

# New row 3
$ws.Range("A3").Value = 77.22916666666666
$ws.Range("B3").Value = 2.0
$ws.Range("C3").Value = 402.0
$ws.Range("D3").Value = 438.0
